$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits: each employee's field block shifts left by one column
#     (the "Naam" field is dropped and a new "Adres4" field is appended
#     after "Adres3"). Apply cell-scoped Replace() calls working from the
#     rightmost affected cell back to the leftmost so that a cell's newly
#     written text is never re-matched by a later Replace() in the same pass.

# Werknemer1 (row 1, columns M..Q)
$ws.Range("Q1").Replace("Adres3Werknemer1", "Adres4Werknemer1")
$ws.Range("P1").Replace("Adres2Werknemer1", "Adres3Werknemer1")
$ws.Range("O1").Replace("Adres1Werknemer1", "Adres2Werknemer1")
$ws.Range("N1").Replace("AdresWerknemer1", "Adres1Werknemer1")
$ws.Range("M1").Replace("NaamWerknemer1", "AdresWerknemer1")

# Werknemer2 (row 2, columns M..Q)
$ws.Range("Q2").Replace("Adres3Werknemer2", "Adres4Werknemer2")
$ws.Range("P2").Replace("Adres2Werknemer2", "Adres3Werknemer2")
$ws.Range("O2").Replace("Adres1Werknemer2", "Adres2Werknemer2")
$ws.Range("N2").Replace("AdresWerknemer2", "Adres1Werknemer2")
$ws.Range("M2").Replace("NaamWerknemer2", "AdresWerknemer2")

# --- Selection moved to A1, with A1:T3 selected
$ws.Range("A1:T3").Select()

# --- Column widths normalised to a single uniform width across columns A:T
$ws.Range("A1:T1").Columns.ColumnWidth = 10.66

# --- Style touch-ups: cells that should fall back to the default,
#     unstyled look (matches the removal of the now-unused custom font
#     that previously backed the "touched" style on these cells)
$ws.Range("N1").Font.Bold = $false
$ws.Range("M2").Font.Bold = $false
$ws.Range("N2").Font.Bold = $false
$ws.Range("O2").Font.Bold = $false
$ws.Range("R2").Font.Bold = $false
